# Updates cryptos list price/volume figures (and a couple of row-content
# swaps) to match the latest scrape, per commit "Updated cryptos list on
# Mon Jul 31 09:11:23 UTC 2023 with GitHub Actions".
#
# For cells whose new value is a plain decimal number (e.g. "243.92"),
# force the cell to Text format first so Excel keeps the original textual
# representation instead of silently re-typing the cell as a Number (which
# would also round-trip floating point weirdness, e.g. 243.92 ->
# 243.91999999999999). The Style is then reset back to "Normal" so no
# stray formatting is left behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.410.76'
$ws.Range("D3").Value = '1.869.55'
$ws.Range("E3").Value = '  -0.35%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '243.92'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.50%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.7058'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.04%  '
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07937'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.11%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3140'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.60%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '24.51'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.94%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07858'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -4.47%  '
$ws.Range("D12").Value = '1.866.96'
$ws.Range("E12").Value = '  -0.90%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '93.85'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.03%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.187'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.18%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.7028'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.38%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.544'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.18%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000008400'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.80%  '
$ws.Range("D18").Value = '29.417.99'
$ws.Range("E18").Value = '  +0.22%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '254.52'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +4.39%  '
$ws.Range("D20").Value = '2.124.00'
$ws.Range("E20").Value = '  -1.23%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.13'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.95%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.000'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.07%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.657'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.57%  '
$ws.Range("E24").Value = '  -0.08%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1558'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.21%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.014'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.39%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '161.02'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.91%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.87'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.77%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.504'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.24%  '
$ws.Range("E30").Value = '  -2.01%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.257'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.11%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.213'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.40%  '
$ws.Range("E33").Value = '  -1.43%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.898'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.20%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7512'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.72%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.178'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.16%  '
$ws.Range("E37").Value = '  +0.87%  '
$ws.Range("B38").Value = 'Maker'
$ws.Range("C38").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D38").Value = '1.287.56'
$ws.Range("E38").Value = '  +2.54%  '
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01883'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.35%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.765'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.40%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8932'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.34%  '
$ws.Range("B42").Value = 'FraxShare'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.022'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -7.01%  '
$ws.Range("B43").Value = 'Quant'
$ws.Range("C43").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '108.85'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.63%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '71.08'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -4.19%  '
$ws.Range("E45").Value = '  -0.03%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00000000129'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -3.37%  '
$ws.Range("D47").Value = '2.024.75'
$ws.Range("E47").Value = '  -1.06%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.801'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.05%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.580'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.13%  '
$ws.Range("E50").Value = '  -0.86%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.4310'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.06%  '
